$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title (A1) to the new slate date
$ws.Range("A1").Value = "NBA, Saturday 2nd Mar 2024"

# Update matchups in column A (rows 2-6)
$ws.Range("A2").Value = "Atlanta Hawks (26-33) vs Brooklyn Nets (23-36)"
$ws.Range("A3").Value = "Utah Jazz (27-33) vs Miami Heat (33-26)"
$ws.Range("A4").Value = "Portland Trail Blazers (16-42) vs Memphis Grizzlies (20-40)"
$ws.Range("A5").Value = "Denver Nuggets (41-19) vs Los Angeles Lakers (34-28)"
$ws.Range("A6").Value = "Houston Rockets (25-34) vs Phoenix Suns (35-24)"

# Update Ballgorithm win-probability picks in column B (rows 2-6)
$ws.Range("B3").Value = "Utah Jazz (62.07%)"
$ws.Range("B4").Value = "Memphis Grizzlies (68.97%)"
$ws.Range("B5").Value = "Denver Nuggets (82.76%)"
$ws.Range("B6").Value = "Houston Rockets (66.67%)"
$ws.Range("B2").Value = "Brooklyn Nets (51.61%)"

# Update ESPN win-probability picks in column C (rows 2-6)
$ws.Range("C2").Value = "Brooklyn Nets (59.2%)"
$ws.Range("C3").Value = "Miami Heat (70.1%)"
$ws.Range("C4").Value = "Memphis Grizzlies (63.3%)"
$ws.Range("C5").Value = "Denver Nuggets (50.3%)"
$ws.Range("C6").Value = "Phoenix Suns (75.3%)"

# Remove the now-unused rows 7-10 (previous slate had 9 games, new slate only has 5)
$ws.Range("A7:C10").EntireRow.Delete()

# Match the author's final view state: entire column D selected
$ws.Range("D1:D1048576").Select()
